$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Data")

# Update A2 value
$ws.Range("A2").Value = 1408798

# Clear the value in A5 but keep its formatting/style
$ws.Range("A5").ClearContents()

# Activate the "Test Data" sheet and update the selection to B14
$ws.Activate()
$ws.Range("B14").Select()
